$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 531
$ws.Range("D2").Value = 588
$ws.Range("E2").Value = 0.47

$ws.Range("C3").Value = 740
$ws.Range("D3").Value = 768
$ws.Range("E3").Value = 0.49

$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").ClearContents()

$ws.Range("C5").Value = 1022
$ws.Range("D5").Value = 1037
$ws.Range("E5").Value = 0.5

$ws.Range("C6").Value = 1386
$ws.Range("D6").Value = 1248
$ws.Range("E6").Value = 0.53
